# re-run RU 1001; without crop
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("B2").Value = 0.678255122017956
$ws.Range("L2").Value = 0.739047667329172

$ws.Range("B3").Value = 0.553746767090015
$ws.Range("D3").Value = 0.622273519140294
$ws.Range("E3").Value = 0.567933697165467
$ws.Range("F3").Value = 0.757304011850663
$ws.Range("G3").Value = 0.53058303552895
$ws.Range("H3").Value = 0.689059343539561
$ws.Range("I3").Value = 0.593195783835747
$ws.Range("J3").Value = 0.623575700142155
$ws.Range("K3").Value = 0.548090076208069
$ws.Range("L3").Value = 0.490485989036895
$ws.Range("M3").Value = 0.846428606188793
$ws.Range("N3").Value = 0.474912333237318

$ws.Range("B4").Value = 0.647210522905015

$ws.Range("B5").Value = 0.671354010141054
$ws.Range("C5").Value = 0.731501666372061
$ws.Range("L5").Value = 0.63442019211072

$ws.Range("B6").Value = 0.684899159160603
$ws.Range("L6").Value = 0.601821932205212

$ws.Range("B7").Value = 0.618798692282585
$ws.Range("L7").Value = 0.540183622108344
